$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value. All target cells are plain text (inline strings) in the
# source workbook, including numeric-looking and percentage-looking values, so each
# write forces a Text number format first (and resets the style afterwards) to stop
# Excel from auto-coercing the string into a Number/Percentage cell.
$updates = @{
    'D2' = '284.83'
    'E2' = '-10.74%'
    'D3' = '40.26'
    'E3' = '-2.07%'
    'D4' = '5.024'
    'E4' = '-3.63%'
    'D5' = '0.07266'
    'E5' = '-5.96%'
    'D6' = '4.291'
    'E6' = '-0.23%'
    'D7' = '1.509'
    'E7' = '-11.34%'
    'D8' = '0.9122'
    'E8' = '-4.17%'
    'D9' = '0.1199'
    'E9' = '-5.34%'
    'D10' = '0.1733'
    'E10' = '-4.96%'
    'D11' = '0.08661'
    'E11' = '-5.42%'
    'D12' = '0.04183'
    'E12' = '-1.13%'
    'D13' = '0.1049'
    'E13' = '-0.36%'
    'D14' = '0.001281'
    'E14' = '0.05%'
    'B15' = 'TigerCash'
    'C15' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D15' = '0.005836'
    'E15' = '-0.76%'
    'B16' = 'LEO'
    'C16' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D16' = '3.397'
    'E16' = '1.26%'
    'B17' = 'BTSEToken'
    'C17' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'D17' = '2.397'
    'E17' = '-1.16%'
    'B18' = 'BitpandaEcosystemToken'
    'C18' = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
    'D18' = '0.3282'
    'E18' = '-2.15%'
    'B19' = 'MCDex'
    'C19' = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
    'D19' = '7.610'
    'E19' = '1.43%'
    'B20' = 'ProBitToken'
    'C20' = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
    'D20' = '0.1345'
    'E20' = '-0.42%'
    'B21' = 'ZBToken'
    'C21' = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
    'D21' = '0.2892'
    'E21' = '3.99%'
    'B22' = 'CoinExToken'
    'C22' = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
    'D22' = '0.03851'
    'E22' = '-4.07%'
    'D23' = '0.001273'
    'E23' = '0.73%'
    'D24' = '0.003760'
    'E24' = '-11.33%'
    'D25' = '0.0001285'
    'E25' = '1.47%'
    'D26' = '0.0003736'
    'E26' = '-95.02%'
    'D38' = '0.02299'
    'E38' = '-9.46%'
    'D39' = '0.04968'
    'E39' = '-7.32%'
    'B40' = 'CEJI'
    'C40' = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
    'D40' = '0.005521'
    'E40' = '179.25%'
    'B41' = 'KickToken'
    'C41' = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
    'D41' = '0.007672'
    'E41' = '-1.28%'
    'D42' = '0.1262'
    'E42' = '-4.22%'
    'D43' = '0.007396'
    'E43' = '0.70%'
    'D44' = '0.007455'
    'E44' = '-1.87%'
    'D45' = '0.3079'
    'E45' = '-10.47%'
    'D46' = '0.00006462'
    'E46' = '-3.73%'
    'D47' = '0.00000000753'
    'E47' = '0.51%'
    'E48' = '15.19%'
    'E49' = '0.18%'
    'D50' = '0.00002108'
    'E50' = '0.51%'
    'D51' = '0.0002007'
    'E51' = '0.51%'
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
    $range.Style = "Normal"
}
